$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.125.08"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "1.823.44"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8703"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.78%  "

$ws.Range("D12").Value = "1.883.51"
$ws.Range("E12").Value = "  +3.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07632"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.336"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.467"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008638"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.85%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "27.440.91"
$ws.Range("E20").Value = "  +0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.215"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("D24").Value = "2.094.90"
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.874"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.083"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.097"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08899"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.959"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7395"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.448"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.17%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.481"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.072"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05243"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01911"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.922"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.162"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5189"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1626"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4829"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.010"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.633"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("E51").Value = "  -0.86%  "
